$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove all existing rows so the shared-string table is rebuilt fresh
# in the order the new values are written below.
$ws.Rows("1:8").Delete()

# Column A: id + numbers
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Column B: Category + values
$ws.Range("B1").Value = "Category"
$ws.Range("B2").Value = "Visual"
$ws.Range("B3").Value = "Thermal"
$ws.Range("B4").Value = "Sound"
$ws.Range("B5").Value = "Chemical"

# Column C: Description + values
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Visual desc"
$ws.Range("C3").Value = "Thermal desc"
$ws.Range("C4").Value = "Sound desc"
$ws.Range("C5").Value = "Chemical desc"

# Re-apply the bold/centered header style (s="2") to the whole header row,
# since deleting rows 1:8 above also wiped the row's prior formatting.
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# Update selection to C1 as in diff
[void]$ws.Range("C1").Select()
